$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 9448
$ws.Cells.Item(2, 5).Value = 7574
$ws.Cells.Item(2, 6).Value = 0.8016511430990686
$ws.Cells.Item(2, 7).Value = 0.7991137370753324
$ws.Cells.Item(2, 8).Value = 0.1028346976498548
$ws.Cells.Item(2, 9).Value = 0.08217661953998734
$ws.Cells.Item(2, 10).Value = 39239770.164845
$ws.Cells.Item(2, 11).Value = 13420220.8417585
$ws.Cells.Item(2, 13).Value = 13420220.8417585
$ws.Cells.Item(2, 14).Value = 52659991.0066035
$ws.Cells.Item(2, 15).Value = 801375021.0472001
$ws.Cells.Item(2, 16).Value = 783675214.0432
$ws.Cells.Item(2, 17).Value = 0.01674649257749708
$ws.Cells.Item(2, 18).Value = 0.01712472284598593

# Row 3
$ws.Cells.Item(3, 4).Value = 9640
$ws.Cells.Item(3, 5).Value = 7551
$ws.Cells.Item(3, 6).Value = 0.783298755186722
$ws.Cells.Item(3, 7).Value = 0.7811918063314711
$ws.Cells.Item(3, 8).Value = 0.1028830618461131
$ws.Cells.Item(3, 9).Value = 0.08037140492447756
$ws.Cells.Item(3, 10).Value = 40538728.96122567
$ws.Cells.Item(3, 11).Value = 13799638.26741243
$ws.Cells.Item(3, 13).Value = 13799638.26741243
$ws.Cells.Item(3, 14).Value = 54338367.2286381
$ws.Cells.Item(3, 15).Value = 836150914.9537281
$ws.Cells.Item(3, 16).Value = 818670738.9296581
$ws.Cells.Item(3, 17).Value = 0.01650376507472469
$ws.Cells.Item(3, 18).Value = 0.01685615182173761

# Row 4
$ws.Cells.Item(4, 4).Value = 9834
$ws.Cells.Item(4, 5).Value = 7543
$ws.Cells.Item(4, 6).Value = 0.7670327435428107
$ws.Cells.Item(4, 7).Value = 0.7651653479407587
$ws.Cells.Item(4, 8).Value = 0.1029497547394936
$ws.Cells.Item(4, 9).Value = 0.07877358490566039
$ws.Cells.Item(4, 10).Value = 42014448.3095379
$ws.Cells.Item(4, 11).Value = 14217426.49631654
$ws.Cells.Item(4, 13).Value = 14217426.49631654
$ws.Cells.Item(4, 14).Value = 56231874.80585443
$ws.Cells.Item(4, 15).Value = 875322312.16536
$ws.Cells.Item(4, 16).Value = 857873364.1594061
$ws.Cells.Item(4, 17).Value = 0.0162425043880644
$ws.Cells.Item(4, 18).Value = 0.01657287321217578

# Row 5
$ws.Cells.Item(5, 4).Value = 10026
$ws.Cells.Item(5, 5).Value = 7534
$ws.Cells.Item(5, 6).Value = 0.7514462397765809
$ws.Cells.Item(5, 7).Value = 0.7493534911478019
$ws.Cells.Item(5, 8).Value = 0.1030222989116007
$ws.Cells.Item(5, 9).Value = 0.07720011935548041
$ws.Cells.Item(5, 10).Value = 43535019.83059579
$ws.Cells.Item(5, 11).Value = 14644484.75463068
$ws.Cells.Item(5, 13).Value = 14644484.75463068
$ws.Cells.Item(5, 14).Value = 58179504.58522647
$ws.Cells.Item(5, 15).Value = 914181328.8136762
$ws.Cells.Item(5, 16).Value = 896695223.350703
$ws.Cells.Item(5, 17).Value = 0.01601923414213095
$ws.Cells.Item(5, 18).Value = 0.01633161900863961

# Row 6
$ws.Cells.Item(6, 4).Value = 10237
$ws.Cells.Item(6, 5).Value = 7523
$ws.Cells.Item(6, 6).Value = 0.7348832665820064
$ws.Cells.Item(6, 7).Value = 0.7336649112541447
$ws.Cells.Item(6, 8).Value = 0.1030931809118703
$ws.Cells.Item(6, 9).Value = 0.07563584942461479
$ws.Cells.Item(6, 10).Value = 45146428.42506469
$ws.Cells.Item(6, 11).Value = 15078985.98419153
$ws.Cells.Item(6, 13).Value = 15078985.98419153
$ws.Cells.Item(6, 14).Value = 60225414.4092562
$ws.Cells.Item(6, 15).Value = 955175919.0364679
$ws.Cells.Item(6, 16).Value = 937584093.1632864
$ws.Cells.Item(6, 17).Value = 0.01578660609388313
$ws.Cells.Item(6, 18).Value = 0.01608280909855989

# Row 7
$ws.Cells.Item(7, 4).Value = 9448
$ws.Cells.Item(7, 5).Value = 7574
$ws.Cells.Item(7, 6).Value = 0.8016511430990686
$ws.Cells.Item(7, 7).Value = 0.7991137370753324
$ws.Cells.Item(7, 8).Value = 0.1028346976498548
$ws.Cells.Item(7, 9).Value = 0.08217661953998734
$ws.Cells.Item(7, 10).Value = 39239770.164845
$ws.Cells.Item(7, 11).Value = 13420220.8417585
$ws.Cells.Item(7, 13).Value = 13420220.8417585
$ws.Cells.Item(7, 14).Value = 52659991.0066035
$ws.Cells.Item(7, 15).Value = 800637542.8872001
$ws.Cells.Item(7, 16).Value = 782937735.8832
$ws.Cells.Item(7, 17).Value = 0.01676191800020205
$ws.Cells.Item(7, 18).Value = 0.01714085325906497

# Row 8
$ws.Cells.Item(8, 4).Value = 9645
$ws.Cells.Item(8, 5).Value = 7551
$ws.Cells.Item(8, 6).Value = 0.7828926905132193
$ws.Cells.Item(8, 7).Value = 0.7811918063314711
$ws.Cells.Item(8, 8).Value = 0.1028830618461131
$ws.Cells.Item(8, 9).Value = 0.08037140492447756
$ws.Cells.Item(8, 10).Value = 40538728.96122567
$ws.Cells.Item(8, 11).Value = 13799638.26741243
$ws.Cells.Item(8, 13).Value = 13799638.26741243
$ws.Cells.Item(8, 14).Value = 54338367.2286381
$ws.Cells.Item(8, 15).Value = 837666957.2402281
$ws.Cells.Item(8, 16).Value = 820186781.216158
$ws.Cells.Item(8, 17).Value = 0.01647389591786768
$ws.Cells.Item(8, 18).Value = 0.01682499472492178

# Row 9
$ws.Cells.Item(9, 4).Value = 9841
$ws.Cells.Item(9, 5).Value = 7543
$ws.Cells.Item(9, 6).Value = 0.766487145615283
$ws.Cells.Item(9, 7).Value = 0.7651653479407587
$ws.Cells.Item(9, 8).Value = 0.1029497547394936
$ws.Cells.Item(9, 9).Value = 0.07877358490566039
$ws.Cells.Item(9, 10).Value = 42014448.3095379
$ws.Cells.Item(9, 11).Value = 14217426.49631654
$ws.Cells.Item(9, 13).Value = 14217426.49631654
$ws.Cells.Item(9, 14).Value = 56231874.80585443
$ws.Cells.Item(9, 15).Value = 875577079.980539
$ws.Cells.Item(9, 16).Value = 858128131.9745851
$ws.Cells.Item(9, 17).Value = 0.01623777828518826
$ws.Cells.Item(9, 18).Value = 0.01656795292749779

# Row 10
$ws.Cells.Item(10, 4).Value = 10035
$ws.Cells.Item(10, 5).Value = 7534
$ws.Cells.Item(10, 6).Value = 0.7507722969606377
$ws.Cells.Item(10, 7).Value = 0.7493534911478019
$ws.Cells.Item(10, 8).Value = 0.1030222989116007
$ws.Cells.Item(10, 9).Value = 0.07720011935548041
$ws.Cells.Item(10, 10).Value = 43535019.83059579
$ws.Cells.Item(10, 11).Value = 14644484.75463068
$ws.Cells.Item(10, 13).Value = 14644484.75463068
$ws.Cells.Item(10, 14).Value = 58179504.58522647
$ws.Cells.Item(10, 15).Value = 913946124.0779107
$ws.Cells.Item(10, 16).Value = 896460018.6149374
$ws.Cells.Item(10, 17).Value = 0.01602335670431958
$ws.Cells.Item(10, 18).Value = 0.01633590394500464

# Row 11
$ws.Cells.Item(11, 4).Value = 10229
$ws.Cells.Item(11, 5).Value = 7523
$ws.Cells.Item(11, 6).Value = 0.7354580115358295
$ws.Cells.Item(11, 7).Value = 0.7336649112541447
$ws.Cells.Item(11, 8).Value = 0.1030931809118703
$ws.Cells.Item(11, 9).Value = 0.07563584942461479
$ws.Cells.Item(11, 10).Value = 45146428.42506469
$ws.Cells.Item(11, 11).Value = 15078985.98419153
$ws.Cells.Item(11, 13).Value = 15078985.98419153
$ws.Cells.Item(11, 14).Value = 60225414.4092562
$ws.Cells.Item(11, 15).Value = 954343579.3721293
$ws.Cells.Item(11, 16).Value = 936751753.4989479
$ws.Cells.Item(11, 17).Value = 0.01580037452980206
$ws.Cells.Item(11, 18).Value = 0.01609709928790484

# Row 12
$ws.Cells.Item(12, 4).Value = 9458
$ws.Cells.Item(12, 5).Value = 7574
$ws.Cells.Item(12, 6).Value = 0.8008035525481074
$ws.Cells.Item(12, 7).Value = 0.7991137370753324
$ws.Cells.Item(12, 8).Value = 0.1028346976498548
$ws.Cells.Item(12, 9).Value = 0.08217661953998734
$ws.Cells.Item(12, 10).Value = 39239770.164845
$ws.Cells.Item(12, 11).Value = 13420220.8417585
$ws.Cells.Item(12, 13).Value = 13420220.8417585
$ws.Cells.Item(12, 14).Value = 52659991.0066035
$ws.Cells.Item(12, 15).Value = 800122294.0972
$ws.Cells.Item(12, 16).Value = 782422487.0932001
$ws.Cells.Item(12, 17).Value = 0.0167727120476013
$ws.Cells.Item(12, 18).Value = 0.01715214102756216

# Row 13
$ws.Cells.Item(13, 4).Value = 9637
$ws.Cells.Item(13, 5).Value = 7551
$ws.Cells.Item(13, 6).Value = 0.7835425962436443
$ws.Cells.Item(13, 7).Value = 0.7811918063314711
$ws.Cells.Item(13, 8).Value = 0.1028830618461131
$ws.Cells.Item(13, 9).Value = 0.08037140492447756
$ws.Cells.Item(13, 10).Value = 40538728.96122567
$ws.Cells.Item(13, 11).Value = 13799638.26741243
$ws.Cells.Item(13, 13).Value = 13799638.26741243
$ws.Cells.Item(13, 14).Value = 54338367.2286381
$ws.Cells.Item(13, 15).Value = 835784763.7231281
$ws.Cells.Item(13, 16).Value = 818304587.6990581
$ws.Cells.Item(13, 17).Value = 0.01651099525425647
$ws.Cells.Item(13, 18).Value = 0.01686369412423168

# Row 14
$ws.Cells.Item(14, 5).Value = 7543
$ws.Cells.Item(14, 6).Value = 0.7665650406504065
$ws.Cells.Item(14, 7).Value = 0.7651653479407587
$ws.Cells.Item(14, 8).Value = 0.1029497547394936
$ws.Cells.Item(14, 9).Value = 0.07877358490566039
$ws.Cells.Item(14, 10).Value = 42014448.3095379
$ws.Cells.Item(14, 11).Value = 14217426.49631654
$ws.Cells.Item(14, 13).Value = 14217426.49631654
$ws.Cells.Item(14, 14).Value = 56231874.80585443
$ws.Cells.Item(14, 15).Value = 874054288.5903099
$ws.Cells.Item(14, 16).Value = 856605340.5843561
$ws.Cells.Item(14, 17).Value = 0.01626606800276291
$ws.Cells.Item(14, 18).Value = 0.01659740585625785

# Row 15
$ws.Cells.Item(15, 4).Value = 10034
$ws.Cells.Item(15, 5).Value = 7534
$ws.Cells.Item(15, 6).Value = 0.7508471197927048
$ws.Cells.Item(15, 7).Value = 0.7493534911478019
$ws.Cells.Item(15, 8).Value = 0.1030222989116007
$ws.Cells.Item(15, 9).Value = 0.07720011935548041
$ws.Cells.Item(15, 10).Value = 43535019.83059579
$ws.Cells.Item(15, 11).Value = 14644484.75463068
$ws.Cells.Item(15, 13).Value = 14644484.75463068
$ws.Cells.Item(15, 14).Value = 58179504.58522647
$ws.Cells.Item(15, 15).Value = 913242019.3379748
$ws.Cells.Item(15, 16).Value = 895755913.8750015
$ws.Cells.Item(15, 17).Value = 0.01603571062712021
$ws.Cells.Item(15, 18).Value = 0.01634874470577512

# Row 16
$ws.Cells.Item(16, 4).Value = 10228
$ws.Cells.Item(16, 5).Value = 7523
$ws.Cells.Item(16, 6).Value = 0.7355299178725069
$ws.Cells.Item(16, 7).Value = 0.7336649112541447
$ws.Cells.Item(16, 8).Value = 0.1030931809118703
$ws.Cells.Item(16, 9).Value = 0.07563584942461479
$ws.Cells.Item(16, 10).Value = 45146428.42506469
$ws.Cells.Item(16, 11).Value = 15078985.98419153
$ws.Cells.Item(16, 13).Value = 15078985.98419153
$ws.Cells.Item(16, 14).Value = 60225414.4092562
$ws.Cells.Item(16, 15).Value = 954929691.5795953
$ws.Cells.Item(16, 16).Value = 937337865.7064139
$ws.Cells.Item(16, 17).Value = 0.01579067665101987
$ws.Cells.Item(16, 18).Value = 0.01608703386033319

Write-Host "Applied all changes"